# "My english vocabulary.xlsx" update — add new words
# Adds 9 new word/translation pairs to the vocabulary list and widens
# column B slightly to fit the longer entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The last existing row (168: "advance" / "продвижение, успех") is
# pushed down to row 169 (a blank row 168 is left behind, matching
# the source workbook), then the new vocabulary rows are appended
# from row 170 through row 178.
# ------------------------------------------------------------------

$ws.Cells.Item(168, 2).Copy($ws.Cells.Item(169, 2))
$ws.Cells.Item(168, 4).Copy($ws.Cells.Item(169, 4))
$ws.Range("B168:D168").Clear()

$ws.Cells.Item(170, 2).Value = "liability"
$ws.Cells.Item(170, 4).Value = "ответственность"

$ws.Cells.Item(171, 2).Value = "do me a favor"
$ws.Cells.Item(172, 2).Value = "you got to be kidding me"
$ws.Cells.Item(172, 4).Value = "ты шутишь надо мной?"
$ws.Cells.Item(171, 4).Value = "сделай мне одолжение, или услугу"

$newWords = @(
    @(173, "rethink", "передумать"),
    @(174, "can I look at it", "могу я посмотреть на это?"),
    @(175, "what are you after", "что ты хочешь,ищешь"),
    @(176, "Wacked", "что-то странное"),
    @(177, "relentless", "неустанный, непрекращающийся"),
    @(178, "map out", "наметить")
)

foreach ($pair in $newWords) {
    $ws.Cells.Item($pair[0], 2).Value = $pair[1]
    $ws.Cells.Item($pair[0], 4).Value = $pair[2]
}

# Column B needs to be a bit wider to fit the new, longer phrases.
$ws.Columns.Item(2).ColumnWidth = 23.95

# Leave the selection where the author left it.
[void]$ws.Range("D167").Select()
